$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83").Value = "2025-04-29 13:21:42"
$ws.Range("B83").Value = 287
